$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2338308457711443
$ws.Range("C2").Value = 0.4577114427860697
$ws.Range("J2").Value = 0.02985074626865672
$ws.Range("P2").Value = 0.1417910447761194
$ws.Range("S2").Value = 0.1368159203980099
$ws.Range("B3").Value = 0.005291005291005291
$ws.Range("C3").Value = 0.02645502645502645
$ws.Range("J3").Value = 0.03703703703703703
$ws.Range("P3").Value = 0.656084656084656
$ws.Range("S3").Value = 0.2751322751322751
$ws.Range("J4").Value = 0.05555555555555555
$ws.Range("P4").Value = 0.5555555555555556
$ws.Range("S4").Value = 0.3888888888888889
$ws.Range("B6").Value = 0.06986899563318777
$ws.Range("D6").Value = 0.008733624454148471
$ws.Range("F6").Value = 0.07423580786026202
$ws.Range("J6").Value = 0.2838427947598253
$ws.Range("O6").Value = 0.02620087336244541
$ws.Range("Q6").Value = 0.1965065502183406
$ws.Range("R6").Value = 0.04803493449781659
$ws.Range("S6").Value = 0.2925764192139738
$ws.Range("B7").Value = 0.08333333333333333
$ws.Range("D7").Value = 0.01785714285714286
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.1845238095238095
$ws.Range("O7").Value = 0.02380952380952381
$ws.Range("Q7").Value = 0.1607142857142857
$ws.Range("R7").Value = 0.04761904761904762
$ws.Range("S7").Value = 0.4345238095238095
$ws.Range("B8").Value = 0.096045197740113
$ws.Range("D8").Value = 0.01318267419962335
$ws.Range("E8").Value = 0.001883239171374765
$ws.Range("F8").Value = 0.04896421845574388
$ws.Range("J8").Value = 0.1525423728813559
$ws.Range("O8").Value = 0.01318267419962335
$ws.Range("Q8").Value = 0.192090395480226
$ws.Range("R8").Value = 0.0903954802259887
$ws.Range("S8").Value = 0.391713747645951
$ws.Range("B9").Value = 0.1058823529411765
$ws.Range("D9").Value = 0.01764705882352941
$ws.Range("F9").Value = 0.03529411764705882
$ws.Range("J9").Value = 0.1764705882352941
$ws.Range("O9").Value = 0.02352941176470588
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.1176470588235294
$ws.Range("S9").Value = 0.3470588235294118
$ws.Range("B10").Value = 0.1270543615676359
$ws.Range("D10").Value = 0.01390644753476612
$ws.Range("E10").Value = 0.0006321112515802782
$ws.Range("F10").Value = 0.06890012642225031
$ws.Range("J10").Value = 0.152338811630847
$ws.Range("O10").Value = 0.01517067003792667
$ws.Range("Q10").Value = 0.2243994943109987
$ws.Range("R10").Value = 0.06005056890012642
$ws.Range("S10").Value = 0.3375474083438685
$ws.Range("G11").Value = 0.1400778210116732
$ws.Range("J11").Value = 0.08171206225680934
$ws.Range("K11").Value = 0.1906614785992218
$ws.Range("L11").Value = 0.5680933852140078
$ws.Range("S11").Value = 0.01945525291828794
$ws.Range("G12").Value = 0.7333333333333333
$ws.Range("J12").Value = 0.2466666666666667
$ws.Range("L12").Value = 0.02
$ws.Range("G13").Value = 0.6341463414634146
$ws.Range("J13").Value = 0.2926829268292683
$ws.Range("S13").Value = 0.07317073170731707
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01298701298701299
$ws.Range("H15").Value = 0.1731601731601732
$ws.Range("I15").Value = 0.05627705627705628
$ws.Range("J15").Value = 0.3766233766233766
$ws.Range("K15").Value = 0.05194805194805195
$ws.Range("M15").Value = 0.008658008658008658
$ws.Range("O15").Value = 0.05627705627705628
$ws.Range("S15").Value = 0.2640692640692641
$ws.Range("F16").Value = 0.02577319587628866
$ws.Range("H16").Value = 0.1701030927835052
$ws.Range("I16").Value = 0.07216494845360824
$ws.Range("J16").Value = 0.4690721649484536
$ws.Range("K16").Value = 0.05154639175257732
$ws.Range("O16").Value = 0.06701030927835051
$ws.Range("S16").Value = 0.1443298969072165
$ws.Range("F17").Value = 0.01444043321299639
$ws.Range("H17").Value = 0.1931407942238267
$ws.Range("I17").Value = 0.0776173285198556
$ws.Range("J17").Value = 0.4476534296028881
$ws.Range("K17").Value = 0.06137184115523465
$ws.Range("M17").Value = 0.01805054151624549
$ws.Range("O17").Value = 0.06859205776173286
$ws.Range("S17").Value = 0.1191335740072202
$ws.Range("F18").Value = 0.01666666666666667
$ws.Range("H18").Value = 0.1722222222222222
$ws.Range("I18").Value = 0.05555555555555555
$ws.Range("J18").Value = 0.5055555555555555
$ws.Range("K18").Value = 0.04444444444444445
$ws.Range("M18").Value = 0.01111111111111111
$ws.Range("O18").Value = 0.04444444444444445
$ws.Range("S18").Value = 0.15
$ws.Range("F19").Value = 0.006419400855920114
$ws.Range("H19").Value = 0.2282453637660485
$ws.Range("I19").Value = 0.06633380884450785
$ws.Range("J19").Value = 0.3873038516405136
$ws.Range("K19").Value = 0.09985734664764621
$ws.Range("M19").Value = 0.02068473609129814
$ws.Range("N19").Value = 0.001426533523537803
$ws.Range("O19").Value = 0.06633380884450785
$ws.Range("S19").Value = 0.12339514978602